# Updates the cryptocurrency price table (columns B-E, rows 2-51) on the
# active worksheet to match the latest scraped values.
#
# Column D ("Price") holds numeric-looking text such as "1.005" or
# "5.290" that must stay exactly as typed (including trailing zeros and
# the "." used as a thousands separator in some rows, e.g. "26.472.72").
# Assigning such a string straight to Range.Value lets Excel silently
# reinterpret it as a number (dropping trailing zeros, etc.), so for
# column D we temporarily switch the cell to Text format, set the
# value, and then restore the cell's original style.

function Set-CellText($ws, $cellRef, $value) {
    $cell = $ws.Range($cellRef)
    $originalStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = $originalStyle
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = "D2"; Value = '26.472.72' },
    @{ Cell = "E2"; Value = '  -0.14%  ' },
    @{ Cell = "D3"; Value = '1.806.71' },
    @{ Cell = "E3"; Value = '  +0.11%  ' },
    @{ Cell = "D4"; Value = '1.005' },
    @{ Cell = "E4"; Value = '  -0.52%  ' },
    @{ Cell = "E5"; Value = '  -0.51%  ' },
    @{ Cell = "D6"; Value = '306.48' },
    @{ Cell = "E6"; Value = '  -0.56%  ' },
    @{ Cell = "D7"; Value = '0.4525' },
    @{ Cell = "E7"; Value = '  -0.49%  ' },
    @{ Cell = "D8"; Value = '0.3599' },
    @{ Cell = "D9"; Value = '46.35' },
    @{ Cell = "E9"; Value = '  +1.29%  ' },
    @{ Cell = "D10"; Value = '0.07072' },
    @{ Cell = "E10"; Value = '  -0.73%  ' },
    @{ Cell = "D11"; Value = '0.8886' },
    @{ Cell = "E11"; Value = '  +2.00%  ' },
    @{ Cell = "D12"; Value = '0.07817' },
    @{ Cell = "E12"; Value = '  +0.44%  ' },
    @{ Cell = "D13"; Value = '19.45' },
    @{ Cell = "E13"; Value = '  +1.37%  ' },
    @{ Cell = "D14"; Value = '1.813.21' },
    @{ Cell = "E14"; Value = '  -1.30%  ' },
    @{ Cell = "D15"; Value = '5.290' },
    @{ Cell = "E15"; Value = '  +0.29%  ' },
    @{ Cell = "D16"; Value = '6.311' },
    @{ Cell = "E16"; Value = '  -0.21%  ' },
    @{ Cell = "D17"; Value = '85.39' },
    @{ Cell = "E17"; Value = '  -1.05%  ' },
    @{ Cell = "E18"; Value = '  -0.43%  ' },
    @{ Cell = "D19"; Value = '0.000008497' },
    @{ Cell = "E19"; Value = '  -0.82%  ' },
    @{ Cell = "E20"; Value = '  -0.47%  ' },
    @{ Cell = "D21"; Value = '26.507.16' },
    @{ Cell = "E21"; Value = '  -0.21%  ' },
    @{ Cell = "E22"; Value = '  -0.20%  ' },
    @{ Cell = "D23"; Value = '4.965' },
    @{ Cell = "E23"; Value = '  +0.16%  ' },
    @{ Cell = "B24"; Value = 'WrappedliquidstakedEther2.0' },
    @{ Cell = "C24"; Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth' },
    @{ Cell = "D24"; Value = '2.045.89' },
    @{ Cell = "E24"; Value = '  -0.44%  ' },
    @{ Cell = "B25"; Value = 'Cosmos' },
    @{ Cell = "C25"; Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom' },
    @{ Cell = "D25"; Value = '10.51' },
    @{ Cell = "E25"; Value = '  +1.51%  ' },
    @{ Cell = "D26"; Value = '1.959' },
    @{ Cell = "E26"; Value = '  -1.34%  ' },
    @{ Cell = "D27"; Value = '151.59' },
    @{ Cell = "E27"; Value = '  +0.25%  ' },
    @{ Cell = "D28"; Value = '17.80' },
    @{ Cell = "E28"; Value = '  -0.24%  ' },
    @{ Cell = "D29"; Value = '2.061' },
    @{ Cell = "E29"; Value = '  +3.08%  ' },
    @{ Cell = "D30"; Value = '112.00' },
    @{ Cell = "E30"; Value = '  -0.69%  ' },
    @{ Cell = "D31"; Value = '4.860' },
    @{ Cell = "E31"; Value = '  -0.11%  ' },
    @{ Cell = "D32"; Value = '0.08690' },
    @{ Cell = "E32"; Value = '  -0.03%  ' },
    @{ Cell = "D33"; Value = '3.099' },
    @{ Cell = "E33"; Value = '  +0.43%  ' },
    @{ Cell = "D34"; Value = '2.849' },
    @{ Cell = "E34"; Value = '  +13.68%  ' },
    @{ Cell = "D35"; Value = '4.458' },
    @{ Cell = "E35"; Value = '  +0.55%  ' },
    @{ Cell = "D36"; Value = '0.7231' },
    @{ Cell = "E36"; Value = '  -0.99%  ' },
    @{ Cell = "D37"; Value = '1.103' },
    @{ Cell = "E37"; Value = '  -0.74%  ' },
    @{ Cell = "D38"; Value = '1.074' },
    @{ Cell = "E38"; Value = '  -0.03%  ' },
    @{ Cell = "D39"; Value = '0.01932' },
    @{ Cell = "E39"; Value = '  +0.96%  ' },
    @{ Cell = "E40"; Value = '  +0.35%  ' },
    @{ Cell = "D41"; Value = '2.892' },
    @{ Cell = "E41"; Value = '  +1.12%  ' },
    @{ Cell = "D42"; Value = '0.5127' },
    @{ Cell = "E42"; Value = '  +4.72%  ' },
    @{ Cell = "D43"; Value = '6.769' },
    @{ Cell = "E43"; Value = '  -1.43%  ' },
    @{ Cell = "D44"; Value = '0.1510' },
    @{ Cell = "E44"; Value = '  -3.54%  ' },
    @{ Cell = "D45"; Value = '8.020' },
    @{ Cell = "E45"; Value = '  -1.26%  ' },
    @{ Cell = "D46"; Value = '0.4664' },
    @{ Cell = "E46"; Value = '  +1.68%  ' },
    @{ Cell = "D47"; Value = '1.003' },
    @{ Cell = "E47"; Value = '  -0.56%  ' },
    @{ Cell = "D48"; Value = '9.973' },
    @{ Cell = "E48"; Value = '  +0.87%  ' },
    @{ Cell = "D49"; Value = '100.19' },
    @{ Cell = "E49"; Value = '  -1.46%  ' },
    @{ Cell = "D50"; Value = '1.573' },
    @{ Cell = "E50"; Value = '  -0.50%  ' },
    @{ Cell = "D51"; Value = '0.05978' },
    @{ Cell = "E51"; Value = '  -0.36%  ' }
)

foreach ($u in $updates) {
    if ($u.Cell.StartsWith("D")) {
        Set-CellText $ws $u.Cell $u.Value
    } else {
        $ws.Range($u.Cell).Value = $u.Value
    }
}
